$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$bsfvbp = $wb.Worksheets.Item("BSfVBP")

# --- About sheet updates ---
$about.Range("B3").Value = "none"
$about.Range("B4").ClearContents()
$about.Range("B5").ClearContents()
$about.Range("B6").ClearContents()
$about.Range("A9").Value = "In the EU only very specific projects receive funding (e.g. via Projects of Common Interest) but there is no general financial support for any battery production."
$about.Range("A10").Value = "That is why no financial support is used here. "
$about.Range("A12").ClearContents()
$about.Range("B12").ClearContents()

$about.Range("A10").Select() | Out-Null

# --- BSfVBP sheet updates ---
$bsfvbp.Range("D2:M2").ClearContents()
$bsfvbp.Range("D2:M2").Value = 0

$bsfvbp.Range("C2").Select() | Out-Null
$bsfvbp.Activate()
